$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.172092999999999
$ws.Range("H2").Value = 21.516279
$ws.Range("I2").Value = 0.07357387076805701
$ws.Range("J2").Value = 0.07357387076805699
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 297.8183156666666
$ws.Range("N2").Value = 893.4549469999999
$ws.Range("O2").Value = 0.8852156413092672
$ws.Range("P2").Value = 0.8852156413092673
$ws.Range("Q2").Value = 2135.98065706469
$ws.Range("R2").Value = 19223.82591358221
$ws.Range("S2").Value = 0.06512874119555073
$ws.Range("T2").Value = 0.06512874119555073

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.172092999999999
$ws.Range("H3").Value = 21.516279
$ws.Range("I3").Value = 0.07357387076805701
$ws.Range("J3").Value = 0.07357387076805699
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 24.34034433333333
$ws.Range("N3").Value = 73.021033
$ws.Range("O3").Value = 0.07234764413494278
$ws.Range("P3").Value = 0.0723476441349428
$ws.Range("Q3").Value = 174.5712132106896
$ws.Range("R3").Value = 1571.140918896207
$ws.Range("S3").Value = 0.005322896219957658
$ws.Range("T3").Value = 0.005322896219957658

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.172092999999999
$ws.Range("H4").Value = 21.516279
$ws.Range("I4").Value = 0.07357387076805701
$ws.Range("J4").Value = 0.07357387076805699
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.277234
$ws.Range("N4").Value = 42.831702
$ws.Range("O4").Value = 0.04243671455578994
$ws.Range("P4").Value = 0.04243671455578994
$ws.Range("Q4").Value = 102.397650030762
$ws.Range("R4").Value = 921.5788502768579
$ws.Range("S4").Value = 0.003122233352548613
$ws.Range("T4").Value = 0.003122233352548612

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 89.72947699999999
$ws.Range("H5").Value = 269.188431
$ws.Range("I5").Value = 0.9204767624852804
$ws.Range("J5").Value = 0.9204767624852804
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 297.8183156666666
$ws.Range("N5").Value = 893.4549469999999
$ws.Range("O5").Value = 0.8852156413092672
$ws.Range("P5").Value = 0.8852156413092673
$ws.Range("Q5").Value = 26723.0817057909
$ws.Range("R5").Value = 240507.7353521181
$ws.Range("S5").Value = 0.8148204276136856
$ws.Range("T5").Value = 0.8148204276136857

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 89.72947699999999
$ws.Range("H6").Value = 269.188431
$ws.Range("I6").Value = 0.9204767624852804
$ws.Range("J6").Value = 0.9204767624852804
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 24.34034433333333
$ws.Range("N6").Value = 73.021033
$ws.Range("O6").Value = 0.07234764413494278
$ws.Range("P6").Value = 0.0723476441349428
$ws.Range("Q6").Value = 2184.046367029914
$ws.Range("R6").Value = 19656.41730326922
$ws.Range("S6").Value = 0.06659432524676932
$ws.Range("T6").Value = 0.06659432524676934

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 89.72947699999999
$ws.Range("H7").Value = 269.188431
$ws.Range("I7").Value = 0.9204767624852804
$ws.Range("J7").Value = 0.9204767624852804
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.277234
$ws.Range("N7").Value = 42.831702
$ws.Range("O7").Value = 0.04243671455578994
$ws.Range("P7").Value = 0.04243671455578994
$ws.Range("Q7").Value = 1281.088739826618
$ws.Range("R7").Value = 11529.79865843956
$ws.Range("S7").Value = 0.03906200962482549
$ws.Range("T7").Value = 0.0390620096248255

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.5799533333333334
$ws.Range("H8").Value = 1.73986
$ws.Range("I8").Value = 0.005949366746662454
$ws.Range("J8").Value = 0.005949366746662453
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 297.8183156666666
$ws.Range("N8").Value = 893.4549469999999
$ws.Range("O8").Value = 0.8852156413092672
$ws.Range("P8").Value = 0.8852156413092673
$ws.Range("Q8").Value = 172.7207248986022
$ws.Range("R8").Value = 1554.48652408742
$ws.Range("S8").Value = 0.005266472500030833
$ws.Range("T8").Value = 0.005266472500030833

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.5799533333333334
$ws.Range("H9").Value = 1.73986
$ws.Range("I9").Value = 0.005949366746662454
$ws.Range("J9").Value = 0.005949366746662453
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 24.34034433333333
$ws.Range("N9").Value = 73.021033
$ws.Range("O9").Value = 0.07234764413494278
$ws.Range("P9").Value = 0.0723476441349428
$ws.Range("Q9").Value = 14.11626383059778
$ws.Range("R9").Value = 127.04637447538
$ws.Range("S9").Value = 0.0004304226682157975
$ws.Range("T9").Value = 0.0004304226682157975

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.5799533333333334
$ws.Range("H10").Value = 1.73986
$ws.Range("I10").Value = 0.005949366746662454
$ws.Range("J10").Value = 0.005949366746662453
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.277234
$ws.Range("N10").Value = 42.831702
$ws.Range("O10").Value = 0.04243671455578994
$ws.Range("P10").Value = 0.04243671455578994
$ws.Range("Q10").Value = 8.280129449080002
$ws.Range("R10").Value = 74.52116504172001
$ws.Range("S10").Value = 0.0002524715784158232
$ws.Range("T10").Value = 0.0002524715784158232
